# Apply the product-stock update described by the diff:
#   - row 6 (Chitos paquete grande): quantity/purchase_price/sale_price were
#     stored as text ("50"/"10000"/"12000"); normalize them to real numbers,
#     and refresh the creation_date timestamp in G6.
#   - add a new row 7 for "Lapicero Retractil kilometrico rojo" (Papelería),
#     whose quantity/purchase_price/sale_price are entered as text (matching
#     how the web form stores brand-new stock rows) while creation_date (G7)
#     is a real datetime number using the same date style as the other rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: normalize quantity / purchase_price / sale_price to numbers ---
$ws.Range("D6").Value = 50
$ws.Range("E6").Value = 10000
$ws.Range("F6").Value = 12000

# --- Row 6: refresh creation_date ---
$ws.Range("G6").Value = 45803.80732170139

# --- Row 7: new stock item ---
$ws.Range("A7").Value = "1P"
$ws.Range("B7").Value = "Papelería"
$ws.Range("C7").Value = "Lapicero Retractil kilometrico rojo"

# quantity / purchase_price / sale_price are kept as plain text (no number
# style lingers on the cells once written, matching the source data).
$newTextCells = $ws.Range("D7:F7")
$newTextCells.NumberFormat = "@"
$ws.Range("D7").Value = "45"
$ws.Range("E7").Value = "1200"
$ws.Range("F7").Value = "2500"
$newTextCells.ClearFormats()

# creation_date: numeric date/time value, formatted like the other rows' G column.
$ws.Range("G7").Value = 45806.85473642316
$ws.Range("G7").NumberFormat = "YYYY-MM-DD HH:MM:SS"
